# Disaggregated.xlsx template maintenance:
#  - drop the unused "Sheet2" stub sheet
#  - extend Sheet1's bordered block from row 37 down to row 110
#  - move the active selection to the new bottom of the sheet

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the empty "Sheet2" placeholder sheet.
$sheet2 = $wb.Worksheets.Item("Sheet2")
[void]$sheet2.Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# Extend the bordered-but-empty block of rows (currently ending at row 37)
# down through row 110 by copying row 37's formatting into the new rows.
[void]$ws.Range("A37:N37").Copy($ws.Range("A38:N110"))

# Scroll/select near the new bottom of the sheet.
[void]$ws.Range("A100").Select()
[void]$ws.Range("N113").Select()
